# Append new sales rows (29-37) to the "Datos" sheet, matching the rows
# produced by the point-of-sale app export (SheetJS) for 12/2/2026.
#
# Columns: A=ID Venta, B=Fecha, C=Hora, D=Vendedor, E=Productos,
#          F=Total (number), G=Personas (number)
#
# Columns A-E must stay literal TEXT (never auto-converted to dates/numbers
# by Excel's type inference), columns F-G are plain numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-breaking space (U+00A0), as used by the app between "p." and "m."
$nbsp = [char]0x00A0

$rows = @(
    @{ Row=29; Id="V-1770933831253"; Fecha="12/2/2026"; Hora="05:03 p.${nbsp}m."; Vendedor="Laura";  Productos="Cerveza Poker (x3)"; Total=27000;  Personas=3 },
    @{ Row=30; Id="V-1770940838425"; Fecha="12/2/2026"; Hora="07:00 p.${nbsp}m."; Vendedor='{"nombre":"Stiven","rol":"ADMIN"}'; Productos="Aguardiente Ligth Caja (x2), Aguardiente Amarillo Botella (x1)"; Total=328000; Personas=0 },
    @{ Row=31; Id="V-1770941236921"; Fecha="12/2/2026"; Hora="07:07 p.${nbsp}m."; Vendedor="Laura";  Productos="Aguardiente Amarillo Caja (x1)"; Total=122000; Personas=2 },
    @{ Row=32; Id="V-1770947115092"; Fecha="12/2/2026"; Hora="08:45 p.${nbsp}m."; Vendedor='{"nombre":"Stiven","rol":"ADMIN"}'; Productos="Aguardiente Amarillo Caja (x1)"; Total=122000; Personas=2 },
    @{ Row=33; Id="V-1770947176403"; Fecha="12/2/2026"; Hora="08:46 p.${nbsp}m."; Vendedor='{"nombre":"Stiven","rol":"ADMIN"}'; Productos="Aguardiente Amarillo Caja (x1)"; Total=122000; Personas=0 },
    @{ Row=34; Id="V-1770947980175"; Fecha="12/2/2026"; Hora="08:59 p.${nbsp}m."; Vendedor="Stiven"; Productos="Aguardiente Ligth Caja (x1), Aguardiente Amarillo Botella (x1), Aguardiente Amarillo Media (x1)"; Total=288000; Personas=0 },
    @{ Row=35; Id="V-1770948065858"; Fecha="12/2/2026"; Hora="09:01 p.${nbsp}m."; Vendedor="Laura";  Productos="Aguardiente Amarillo Media (x1)"; Total=70000;  Personas=2 },
    @{ Row=36; Id="V-1770948134526"; Fecha="12/2/2026"; Hora="09:02 p.${nbsp}m."; Vendedor="Juanita"; Productos="Aguardiente Amarillo Caja (x1)"; Total=122000; Personas=3 },
    @{ Row=37; Id="V-1770950032583"; Fecha="12/2/2026"; Hora="09:33 p.${nbsp}m."; Vendedor="Juanita"; Productos="Aguardiente Amarillo Caja (x1)"; Total=122000; Personas=0 }
)

foreach ($r in $rows) {
    $n = $r.Row

    $ws.Range("A$n").Value = $r.Id

    # Fecha / Hora look numeric-ish (dates/times) to Excel's auto-detection,
    # so force the cell to Text format *before* assigning the value to keep
    # them as literal strings instead of being coerced into date/time serials.
    $ws.Range("B$n").NumberFormat = "@"
    $ws.Range("B$n").Value = $r.Fecha
    $ws.Range("B$n").Style = "Normal"

    $ws.Range("C$n").NumberFormat = "@"
    $ws.Range("C$n").Value = $r.Hora
    $ws.Range("C$n").Style = "Normal"

    $ws.Range("D$n").Value = $r.Vendedor
    $ws.Range("E$n").Value = $r.Productos

    $ws.Range("F$n").Value = $r.Total
    $ws.Range("G$n").Value = $r.Personas
}
